$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the value in A2 from 3 to 5
$ws.Range("A2").Value = 5

# Move the active selection to A3 (cell below the data, no value)
$ws.Range("A3").Select()
